$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '30.465.80'
$ws.Range("E2").Value = '  +0.09%  '
Set-TextCell $ws.Range("D3") '1.913.28'
$ws.Range("E3").Value = '  -0.16%  '
Set-TextCell $ws.Range("D4") '0.9984'
$ws.Range("E4").Value = '  -0.16%  '
Set-TextCell $ws.Range("D5") '244.65'
$ws.Range("E5").Value = '  +0.39%  '
Set-TextCell $ws.Range("D6") '0.9976'
$ws.Range("E6").Value = '  -0.24%  '
Set-TextCell $ws.Range("D7") '0.4787'
$ws.Range("E7").Value = '  +2.00%  '
$ws.Range("E8").Value = '  +0.71%  '
Set-TextCell $ws.Range("D9") '0.06726'
$ws.Range("E9").Value = '  -1.62%  '
Set-TextCell $ws.Range("D10") '111.79'
$ws.Range("E10").Value = '  +1.28%  '
Set-TextCell $ws.Range("D11") '19.27'
$ws.Range("E11").Value = '  +4.71%  '
Set-TextCell $ws.Range("D12") '1.906.03'
$ws.Range("E12").Value = '  -0.47%  '
Set-TextCell $ws.Range("D13") '0.07550'
$ws.Range("E13").Value = '  -2.45%  '
Set-TextCell $ws.Range("D14") '5.237'
$ws.Range("E14").Value = '  -1.05%  '
Set-TextCell $ws.Range("D15") '0.6667'
$ws.Range("E15").Value = '  +1.48%  '
Set-TextCell $ws.Range("D16") '303.79'
$ws.Range("E16").Value = '  +2.75%  '
Set-TextCell $ws.Range("D17") '30.469.20'
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("E18").Value = '  -0.06%  '
Set-TextCell $ws.Range("D19") '0.9980'
$ws.Range("E19").Value = '  -0.15%  '
Set-TextCell $ws.Range("D20") '0.000007566'
$ws.Range("E20").Value = '  -0.98%  '
Set-TextCell $ws.Range("D21") '2.161.06'
$ws.Range("E21").Value = '  +0.72%  '
Set-TextCell $ws.Range("D22") '5.473'
$ws.Range("E22").Value = '  +4.35%  '
Set-TextCell $ws.Range("D23") '0.9960'
$ws.Range("E23").Value = '  -0.44%  '
Set-TextCell $ws.Range("D24") '6.397'
$ws.Range("E24").Value = '  +2.96%  '
Set-TextCell $ws.Range("D25") '9.473'
$ws.Range("E25").Value = '  +1.13%  '
Set-TextCell $ws.Range("D26") '164.19'
$ws.Range("E26").Value = '  -2.77%  '
Set-TextCell $ws.Range("D27") '20.61'
$ws.Range("E27").Value = '  -5.08%  '
Set-TextCell $ws.Range("D28") '2.091'
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("E29").Value = '  +0.35%  '
Set-TextCell $ws.Range("D30") '1.399'
$ws.Range("E30").Value = '  +2.58%  '
Set-TextCell $ws.Range("D31") '4.167'
$ws.Range("E31").Value = '  -0.22%  '
Set-TextCell $ws.Range("D32") '4.024'
$ws.Range("E32").Value = '  +1.00%  '
Set-TextCell $ws.Range("D33") '0.04970'
$ws.Range("E33").Value = '  -1.47%  '
Set-TextCell $ws.Range("D34") '0.7330'
$ws.Range("E34").Value = '  -0.34%  '
Set-TextCell $ws.Range("D35") '1.137'
$ws.Range("E35").Value = '  -1.49%  '
Set-TextCell $ws.Range("D36") '0.9988'
$ws.Range("E36").Value = '  +0.02%  '
Set-TextCell $ws.Range("D37") '0.02048'
$ws.Range("E37").Value = '  -1.04%  '
Set-TextCell $ws.Range("D38") '2.723'
$ws.Range("E38").Value = '  -0.65%  '
Set-TextCell $ws.Range("D39") '2.669'
$ws.Range("E39").Value = '  -0.58%  '
Set-TextCell $ws.Range("D40") '111.44'
$ws.Range("E40").Value = '  +1.75%  '
Set-TextCell $ws.Range("D41") '2.021'
$ws.Range("E41").Value = '  -1.82%  '
Set-TextCell $ws.Range("D42") '0.4413'
$ws.Range("E42").Value = '  +3.77%  '
Set-TextCell $ws.Range("D43") '0.8630'
$ws.Range("E43").Value = '  -0.88%  '
Set-TextCell $ws.Range("D44") '5.905'
$ws.Range("E44").Value = '  +1.04%  '
Set-TextCell $ws.Range("D45") '0.9971'
$ws.Range("E45").Value = '  -0.26%  '
Set-TextCell $ws.Range("D46") '68.69'
$ws.Range("E46").Value = '  +1.89%  '
Set-TextCell $ws.Range("D47") '49.82'
$ws.Range("E47").Value = '  -2.98%  '
Set-TextCell $ws.Range("D48") '7.290'
$ws.Range("E48").Value = '  +1.35%  '
Set-TextCell $ws.Range("D49") '9.284'
$ws.Range("E49").Value = '  +0.62%  '
Set-TextCell $ws.Range("D50") '0.1233'
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("E51").Value = '  +4.23%  '
